# Insert a new weekly price record as row 141 in the "Ciruela" sheet,
# shifting the existing rows 141-171 down to 142-172.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 141 (pushes old row 141.. down by one)
$ws.Rows.Item(141).Insert()

# Populate the new row 141 with the new data record
$ws.Cells.Item(141, 1).Value = 4
$ws.Cells.Item(141, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(141, 3).Value = "Los Lagos"
$ws.Cells.Item(141, 4).Value = 44627
$ws.Cells.Item(141, 5).Value = 10
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100103
$ws.Cells.Item(141, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(141, 9).Value = 100103002
$ws.Cells.Item(141, 10).Value = "Ciruela"
$ws.Cells.Item(141, 11).Value = "Black Amber"
$ws.Cells.Item(141, 12).Value = "Primera"
$ws.Cells.Item(141, 13).Value = 300
$ws.Cells.Item(141, 14).Value = 15000
$ws.Cells.Item(141, 15).Value = 15000
$ws.Cells.Item(141, 16).Value = 15000
$ws.Cells.Item(141, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(141, 18).Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Cells.Item(141, 19).Value = 1000
$ws.Cells.Item(141, 20).Value = 15
